$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This edit rotates several chunks of text among existing paragraphs while
# leaving paragraph formatting / styles / run formatting (bold, italic, ...)
# untouched. Capture every source text value BEFORE any paragraph is
# overwritten, then write the captured values into their destinations.
# ---------------------------------------------------------------------------

# Paragraph index map (1-based, as exposed by $d.Paragraphs):
#  6  -> "Objetivos" body (PT, plain)
#  7  -> "Objetivos" body (EN, italic)
#  9  -> "Docente(s) Responsável(eis)" bullet
# 11  -> "Programa resumido" body (PT, plain)
# 12  -> "Programa resumido" body (EN, italic)
# 14  -> "Programa" body (PT, plain)
# 17  -> "Avaliação" bullet paragraph with 3 bold labels + 3 answer runs
# 19  -> "Bibliografia" body

function Get-ParaText($index) {
    $p = $d.Paragraphs.Item($index)
    $t = $p.Range.Text
    # Paragraph.Range.Text includes the trailing paragraph-mark character
    # (CR, 0x0D) - strip it so callers get just the visible content.
    if ($t.Length -gt 0 -and [int][char]$t[$t.Length - 1] -eq 13) {
        $t = $t.Substring(0, $t.Length - 1)
    }
    return $t
}

function Set-ParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    # Leave the trailing paragraph mark (and any trailing line breaks) alone;
    # only replace the visible character content.
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

# --- capture originals ------------------------------------------------------
$orig6  = Get-ParaText 6    # Objetivos PT
$orig7  = Get-ParaText 7    # Objetivos EN
$orig9  = Get-ParaText 9    # Docente bullet
$orig11 = Get-ParaText 11   # Programa resumido PT
$orig12 = Get-ParaText 12   # Programa resumido EN
$orig14 = Get-ParaText 14   # Programa PT

$metodoAnswer  = "Aulas expositivas, seminários e exercícios comentados."
$criterioAnswer = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$normaAnswer = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

$orig19 = Get-ParaText 19    # Bibliografia

# --- write new values --------------------------------------------------------

# "Objetivos" section now holds what used to be "Programa resumido" text
Set-ParaText 6 $orig11
Set-ParaText 7 $orig12

# "Docente(s) Responsável(eis)" bullet now holds what used to be the
# "Objetivos" PT paragraph text
Set-ParaText 9 $orig6

# "Programa resumido" now holds what used to be the "Programa" (long) text
Set-ParaText 11 $orig14

# "Programa resumido" EN (italic) now holds what used to be "Objetivos" EN text
Set-ParaText 12 $orig7

# "Programa" body now holds the old "Método:" answer text
Set-ParaText 14 $metodoAnswer

# Inside the "Avaliação" bullet paragraph, rotate the three answer runs:
# Método: gets old Critério answer; Critério: gets old Norma answer;
# Norma de recuperação: gets the old Bibliografia text.
#
# Apply the replacements right-to-left (Norma, then Critério, then Método).
# Each step's search text is still the untouched original at that point,
# so there is no risk of a later step accidentally re-matching text that
# an earlier step just wrote (which would happen if processed left-to-right,
# since e.g. the new Método text would equal the original Critério text).
$avalRange = $d.Paragraphs.Item(17).Range
$avalRange.Find.Execute($normaAnswer, $true, $false, $false, $false, $false, $true, 1, $false, $orig19, 2) | Out-Null

$avalRange = $d.Paragraphs.Item(17).Range
$avalRange.Find.Execute($criterioAnswer, $true, $false, $false, $false, $false, $true, 1, $false, $normaAnswer, 2) | Out-Null

$avalRange = $d.Paragraphs.Item(17).Range
$avalRange.Find.Execute($metodoAnswer, $true, $false, $false, $false, $false, $true, 1, $false, $criterioAnswer, 2) | Out-Null

# "Bibliografia" body now holds what used to be the "Docente(s)" bullet text
Set-ParaText 19 $orig9
